# Applies the "moduly vs. části" wording updates plus the two content
# additions described in the commit's diff.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $old"
    }
}

# 1) "pro jednotlivé části" -> "pro jednotlivé moduly"
Replace-Text `
    "pro jednotlivé části, které jsou přihlášenému uživateli dostupné dle jeho" `
    "pro jednotlivé moduly, které jsou přihlášenému uživateli dostupné dle jeho"

# 2) "Jednotlivé části AMČR ... konkrétní části" -> "Jednotlivé moduly AMČR ... konkrétního modulu"
Replace-Text `
    "Jednotlivé části AMČR jsou barevně odlišeny a daná barevnost je pak dodržena i při vstupu do konkrétní části, na první pohled je tedy zřejmé, kde se uživatel nachází." `
    "Jednotlivé moduly AMČR jsou barevně odlišeny a daná barevnost je pak dodržena i při vstupu do konkrétního modulu, na první pohled je tedy zřejmé, kde se uživatel nachází."

# 3) Add "dostupná pro uživatele s oprávněním archeolog a výše, " to the "Naše záznamy" filter description
Replace-Text `
    "– Přednastavená filtrace, která zobrazuje záznamy vytvořené všemi uživateli stejné organizace jako je organizace přihlášeného uživatele." `
    "– Přednastavená filtrace, dostupná pro uživatele s oprávněním archeolog a výše, která zobrazuje záznamy vytvořené všemi uživateli stejné organizace jako je organizace přihlášeného uživatele."

# 4) "dle dané části AMČR" -> "dle daného modulu AMČR"
Replace-Text `
    "Aktuální procesní stav je vždy zvýrazněn barvou dle dané části AMČR a pod ukazateli procesních stavů, kterými již záznam prošel, je datum, kdy k poslední změně daného stavu došlo." `
    "Aktuální procesní stav je vždy zvýrazněn barvou dle daného modulu AMČR a pod ukazateli procesních stavů, kterými již záznam prošel, je datum, kdy k poslední změně daného stavu došlo."

# 5) Rewrite the "Po povolení úprav" sentence, then append two new runs
#    after it: a single-space run and a new sentence about the
#    yellow-bordered fields.
$rng = $d.Content
$found = $rng.Find.Execute( `
    "Po povolení úprav se zobrazí konkrétní upravovaná sekce daného záznamu, kdy pole ohraničená žlutě je v daném procesním stavu povinné vyplnit.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Po povolení úprav se zobrazí konkrétní upravovaná sekce daného záznamu, kdy pole označená za popiskem hvězdičkou, jsou povinná.", `
    2)
if (-not $found) {
    throw "Find/Replace failed for the 'Po povolení úprav' sentence"
}

# After Find.Execute replaces the text, $rng spans exactly the new
# sentence. Collapse it to its end and insert the two extra runs there.
$rng.Collapse(0)
$rng.InsertAfter(" ")
$rng.Collapse(0)
$rng.InsertAfter("Pole ohraničená žlutě je pak nutné vyplnit pro posun do dalšího procesního stavu.")
